$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-29 12:16:59"
$wsOverview.Range("G5").Value = "2016-08-29 12:16:59"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-29 12:16:55"
$wsZhCn.Range("H5").Value = "2016-08-29 12:16:55"
$wsZhCn.Range("K2").Value = "2016-08-29 12:17:17"
$wsZhCn.Range("K5").Value = "2016-08-29 12:17:17"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-29 12:16:59"
$wsDeDe.Range("H5").Value = "2016-08-29 12:16:59"
$wsDeDe.Range("K2").Value = "2016-08-29 12:17:24"
$wsDeDe.Range("K5").Value = "2016-08-29 12:17:24"
